$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes --------------------------------------------------
# Insert 3 new rows at position 4 for 2 "under_review" + 1 "in_progress"
# entries; this pushes the existing "conference" rows (old 4-7) down to 7-10.
$ws.Rows("4:6").Insert()

# The row insert copies formatting from the row above into the unused
# columns F:I and K; drop that leftover formatting so the new rows only
# carry the columns that actually have data.
$ws.Range("F4:I6").Clear()
$ws.Range("K4:K6").Clear()

# Insert 2 new columns at L:M for "submission"/"department" fields,
# shifting the existing L..U columns (book..version) right to N..W.
$ws.Range("L1:M11").Insert(-4161)

# --- New header cells ------------------------------------------------------
$ws.Range("L1").Value = 'submission'
$ws.Range("M1").Value = 'department'

# --- New data rows (2x under_review, 1x in_progress) -----------------------
# Row 4
$ws.Range("A4").Value = 'under_review'
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 'Wright, J., Sohlberg, M.M., McIntosh, K., Seeley, J., Hadley, W., Blitz, D. & Lowham, E.'
$ws.Range("D4").Value = 2021
$ws.Range("E4").Value = 'What is the effect of personalized cognitive strategy instruction on facilitating return-to-learn for individuals experiencing prolonged concussion symptoms?'
$ws.Range("L4").Value = '[Manuscript submitted for publication]'
$ws.Range("M4").Value = 'Department of Communication Disorders & Sciences, University of Oregon'

# Row 5
$ws.Range("A5").Value = 'under_review'
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 'Kucheria, P., Sohlberg, M.M., Fickas, S. Prideaux, J., & Wright, J.'
$ws.Range("D5").Value = 2021
$ws.Range("E5").Value = '"RULE"-ing out comprehension deficits: Validity of the RULE tool as a screener for measuring postsecondary reading comprehension '
$ws.Range("L5").Value = '[Manuscript submitted for publication]'
$ws.Range("M5").Value = 'Department of Communication Disorders & Sciences, University of Oregon'

# Row 6
$ws.Range("A6").Value = 'in_progress'
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 'Wright, J., Furutani, T., Sohlberg, M.M., Mashima, P., & Murata, N.'
$ws.Range("D6").Value = 2021
$ws.Range("E6").Value = 'A retrospective review of thirteen years of concussion symptom reporting and trajectory data across the State of Hawaii and its influence on the future of return-to-learn '
$ws.Range("L6").Value = '[Unpublished manuscript]'
$ws.Range("M6").Value = 'Department of Communication Disorders & Scienes, University of Oregon; Department of Kinesiology and Rehabilitation Sciences, University of Hawaii at Manoa; Department of Communication Sciences and Disorders, University of Hawaii at Manoa'

# --- New poster/conference row (row 11) -------------------------------------
$ws.Range("A11").Value = 'conference'
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 'Wright, J. & Loucks, T.'
$ws.Range("D11").Value = '2014, February '
$ws.Range("E11").Value = 'Adaptation under Altered Auditory Feedback'
$ws.Range("P11").Value = 'Illinois Speech-Language-Hearing Association'
$ws.Range("Q11").Value = $false
$ws.Range("R11").Value = 'Rosemont, IL'
$ws.Range("S11").Value = $true

# --- Row heights (match authored/auto-fit heights from the source file) ----
$ws.Rows("4").RowHeight = 119
$ws.Rows("5").RowHeight = 119
$ws.Rows("6").RowHeight = 221
$ws.Rows("7").RowHeight = 136
$ws.Rows("8").RowHeight = 119
$ws.Rows("9").RowHeight = 119
$ws.Rows("10").RowHeight = 136
$ws.Rows("11").RowHeight = 51

# --- View/selection state ----------------------------------------------------
$ws.Range("M7").Select() | Out-Null

